$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.383.19"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "2.061.01"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.03%  "
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D13").Value = "2.364.31"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.775"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "2.060.26"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "37.298.37"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.37%  "
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("E28").Value = "  +5.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  +6.38%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.12%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "1.480.53"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.97%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0932"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
